$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "ბორჯომი"

# Remove the subtitle row "(მოსახლეობის აღწერის შედეგებით)" (old row 2)
$ws.Rows("2").Delete()

# Remove the extra year columns (old C:D, which held the 2002 and 2014 data,
# duplicated with identical area values); the remaining column B will hold
# the single remaining year.
$ws.Range("C:D").Delete()

# The cell left behind in column B (row 4) previously held "1989"; update it
# to the retained year, 2014.
$ws.Range("B4").Value = 2014

# Clear the now-stray, empty formatted cells left behind in column B of the
# title/subtitle rows so they don't linger with dead styles.
$ws.Range("B1").Clear()
$ws.Range("B2").Clear()

# Restore the active cell/selection
$ws.Range("A2").Select()
